$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = "Add -  Edit and Delete Position Category From Excel"

$ws.Range("A2").Value = $title
$ws.Range("B2").Value = "FAILED"
$ws.Range("C2").Value = "chrome"

$ws.Range("A3").Value = $title
$ws.Range("B3").Value = "PASSED"
$ws.Range("C3").Value = "chrome"
